$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values on rows 2-3 ---
$ws.Cells.Item(2,7).Value = 2.99      # G2: Bundle Diameter
$ws.Cells.Item(2,8).Value = 5.25      # H2: Bundle Weight
$ws.Cells.Item(3,5).Value = "7C#14"   # E3: Cable Size

# --- Add new rows 4-10 with Pull # (D) and Cable Size (E) ---
$pulls = @(3,4,5,6,7,8,9)
$sizes = @("7C#14","7C#14","7C#14","7C#14","7C#14","2C#2","2C#2")

for ($i = 0; $i -lt $pulls.Length; $i++) {
    $r = $i + 4
    $ws.Cells.Item($r,4).Value = $pulls[$i]
    $ws.Cells.Item($r,5).Value = $sizes[$i]
}

# --- Apply the same centered style used by rows 2-3 to all new row cells (A:J) ---
for ($r = 4; $r -le 10; $r++) {
    $rowRange = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,10))
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108
}

# --- Re-merge columns A,B,C,F,G,H,I,J across rows 2-10 (un-merging old 2:3 first) ---
$ws.Range("A2:A3").UnMerge()
$ws.Range("B2:B3").UnMerge()
$ws.Range("C2:C3").UnMerge()
$ws.Range("F2:F3").UnMerge()
$ws.Range("G2:G3").UnMerge()
$ws.Range("H2:H3").UnMerge()
$ws.Range("I2:I3").UnMerge()
$ws.Range("J2:J3").UnMerge()

$ws.Range("A2:A10").Merge()
$ws.Range("B2:B10").Merge()
$ws.Range("C2:C10").Merge()
$ws.Range("F2:F10").Merge()
$ws.Range("G2:G10").Merge()
$ws.Range("H2:H10").Merge()
$ws.Range("I2:I10").Merge()
$ws.Range("J2:J10").Merge()
